$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.171.69"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "1.812.20"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.67"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.555"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.03"
$ws.Range("E8").Value = "  -2.87%  "
$ws.Range("E9").Value = "  +3.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0726"
$ws.Range("E10").Value = "  +9.85%  "
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "2.076.14"
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("D13").Value = "1.810.82"
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.97"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.639"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").Value = "34.188.05"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.31"
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.36"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.80"
$ws.Range("E19").Value = "  -2.60%  "
$ws.Range("D20").Value = "0.0₃0796"
$ws.Range("E20").Value = "  +6.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.00"
$ws.Range("E21").Value = "  +5.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.24"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.08"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.66"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.18"
$ws.Range("E27").Value = "  +2.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.114"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0532"
$ws.Range("E30").Value = "  +3.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.77"
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.59"
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.88"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").Value = "1.430.48"
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.638"
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0189"
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.957"
$ws.Range("E39").Value = "  +7.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "81.25"
$ws.Range("E40").Value = "  -3.15%  "
$ws.Range("E41").Value = "  -3.81%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.15"
$ws.Range("E43").Value = "  +3.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.98"
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("D46").Value = "1.971.36"
$ws.Range("E46").Value = "  +1.12%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.04"
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.42"
$ws.Range("E48").Value = "  +7.41%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.95"
$ws.Range("E49").Value = "  -2.40%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.996"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("E51").Value = "  +6.03%  "
